$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new weekly data point (week of 2021-09-10, serial 44449) is inserted
# ahead of the existing "Verde" rows, pushing everything from the old
# row 25 down by two rows (old 25-30 -> new 27-32).
$ws.Rows("25:26").Insert()

# New row 25: Espárragos, "Sin especificar" variety, "Primera" quality
$ws.Range("A25").Value = 9
$ws.Range("B25").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C25").Value = "Metropolitana"
$ws.Range("D25").Value = 44449
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = 300000000
$ws.Range("G25").Value = "Espárragos"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 16
$ws.Range("K25").Value = 34000
$ws.Range("L25").Value = 34000
$ws.Range("M25").Value = 34000
$ws.Range("N25").Value = "$/bandeja 10 kilos"
$ws.Range("O25").Value = "Región Metropolitana"
$ws.Range("P25").Value = 3400
$ws.Range("Q25").Value = 10
$ws.Range("R25").Value = "Hortaliza"

# New row 26: Espárragos, "Sin especificar" variety, "Segunda" quality
$ws.Range("A26").Value = 9
$ws.Range("B26").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44449
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = 300000000
$ws.Range("G26").Value = "Espárragos"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 7
$ws.Range("K26").Value = 27000
$ws.Range("L26").Value = 27000
$ws.Range("M26").Value = 27000
$ws.Range("N26").Value = "$/bandeja 10 kilos"
$ws.Range("O26").Value = "Región Metropolitana"
$ws.Range("P26").Value = 2700
$ws.Range("Q26").Value = 10
$ws.Range("R26").Value = "Hortaliza"
